$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D header "Single Threaded" (bold) — adds a shared string +
# a bold font / cell style, and extends the sheet's used range to D.
$ws.Range("D1").Value = "Single Threaded"
$ws.Range("D1").Font.Bold = $true

# Widen column D to fit the new header text.
$ws.Columns("D").ColumnWidth = 17.43

# Page setup: portrait, paper size 9 (A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to G1 (matches the saved cursor position).
$ws.Range("G1").Select()
